$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FREQ")
$ws2 = $wb.Worksheets.Item("PONTOS")

$values = @{
    3  = "F"
    4  = "F"
    5  = "F"
    6  = "P"
    7  = "F"
    8  = "P"
    9  = "P"
    10 = "P"
    11 = "F"
    12 = "F"
    13 = "P"
    14 = "F"
    15 = "P"
    16 = "P"
    17 = "F"
    18 = "F"
    19 = "P"
    20 = "F"
}

foreach ($row in $values.Keys) {
    $ws1.Range("R$row").Value = $values[$row]
}

$ws1.Range("R16").Select()

$ws2.Range("E5").Value = 2
$ws2.Range("G11").Select()
